{"js": "// This document contains a date/weekday heading paragraph followed by a\n// 20x5 table of simple arithmetic drills (one equation per cell). The edit\n// updates the heading date and rewrites every equation in the table, in\n// document order. Each pair below is [oldText, newText] taken straight\n// from the authoritative OOXML diff, and every oldText is unique in the\n// document, so we can safely match paragraphs by their current text.\nconst REPLACEMENTS = [[\"2023-02-23 Thursday\", \"2023-02-24 Friday\"], [\"8+54=\", \"48-16=\"], [\"42+28=\", \"32+38=\"], [\"67-63=\", \"17+2=\"], [\"8+80=\", \"84-73=\"], [\"6+66=\", \"63+34=\"], [\"24+5=\", \"88-81=\"], [\"41-8=\", \"66-19=\"], [\"39+13=\", \"46-12=\"], [\"97-94=\", \"93-29=\"], [\"39+15=\", \"4+67=\"], [\"95-71=\", \"96-57=\"], [\"53-47=\", \"1+92=\"], [\"94-22=\", \"91-28=\"], [\"0+69=\", \"11+30=\"], [\"90-22=\", \"67+11=\"], [\"32+6=\", \"63-2=\"], [\"10+53=\", \"91-90=\"], [\"29+17=\", \"84-40=\"], [\"73-18=\", \"75-29=\"], [\"52+34=\", \"33+8=\"], [\"89-16=\", \"15+82=\"], [\"61-9=\", \"43+3=\"], [\"9+0=\", \"20+40=\"], [\"64-62=\", \"60+9=\"], [\"80-48=\", \"45-40=\"], [\"96-74=\", \"60-36=\"], [\"17+65=\", \"91-8=\"], [\"13+27=\", \"69-67=\"], [\"33+29=\", \"75+9=\"], [\"28+3=\", \"79-30=\"], [\"13+46=\", \"13+32=\"], [\"34-23=\", \"74+10=\"], [\"91-30=\", \"64-36=\"], [\"77-44=\", \"66-3=\"], [\"39+11=\", \"36-23=\"], [\"41-3=\", \"43-28=\"], [\"84-25=\", \"16+73=\"], [\"89-39=\", \"48-35=\"], [\"70-4=\", \"24+21=\"], [\"76-15=\", \"65+7=\"], [\"1+85=\", \"36-31=\"], [\"80-5=\", \"29+27=\"], [\"37+49=\", \"37+19=\"], [\"23+4=\", \"75+13=\"], [\"70-18=\", \"2+27=\"], [\"34+25=\", \"85-1=\"], [\"34+46=\", \"65+13=\"], [\"49+9=\", \"1+59=\"], [\"55-12=\", \"5+36=\"], [\"35-27=\", \"62-32=\"], [\"65+29=\", \"91-19=\"], [\"19+48=\", \"35+6=\"], [\"61-0=\", \"0+3=\"], [\"58+25=\", \"20+76=\"], [\"44+18=\", \"97-25=\"], [\"35-25=\", \"70-50=\"], [\"36+18=\", \"46-46=\"], [\"74-25=\", \"45+44=\"], [\"37+37=\", \"16+16=\"], [\"67-52=\", \"57+22=\"], [\"76-37=\", \"70+22=\"], [\"32+34=\", \"71+9=\"], [\"4+77=\", \"92-18=\"], [\"25+52=\", \"66-29=\"], [\"65-52=\", \"32+11=\"], [\"66+4=\", \"66+27=\"], [\"94-8=\", \"36+56=\"], [\"5+58=\", \"44+1=\"], [\"44+38=\", \"39+5=\"], [\"70-39=\", \"96-57=\"], [\"27-3=\", \"55+28=\"], [\"41+1=\", \"89-86=\"], [\"86-48=\", \"34+31=\"], [\"42+12=\", \"56+12=\"], [\"6+49=\", \"48-27=\"], [\"95-83=\", \"14+60=\"], [\"33+2=\", \"22+58=\"], [\"85-43=\", \"32+57=\"], [\"93-90=\", \"11+38=\"], [\"96-62=\", \"49-38=\"], [\"65-53=\", \"93-1=\"], [\"28+42=\", \"98-63=\"], [\"87-16=\", \"31+19=\"], [\"40+36=\", \"50+35=\"], [\"10+55=\", \"69-20=\"], [\"32+0=\", \"38+7=\"], [\"45-22=\", \"36-3=\"], [\"32+56=\", \"28+32=\"], [\"65+8=\", \"52-48=\"], [\"7+79=\", \"88-40=\"], [\"66+26=\", \"16+35=\"], [\"36-20=\", \"71+22=\"], [\"52-13=\", \"62-36=\"], [\"53+32=\", \"50+7=\"], [\"0+19=\", \"63-16=\"], [\"24-6=\", \"63+13=\"], [\"34+29=\", \"35+18=\"], [\"43-42=\", \"33+61=\"], [\"10+60=\", \"86-29=\"], [\"15+38=\", \"86-42=\"]];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} paragraphs (1 heading + 100 table cells), found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = paragraphs.items[i];\n  const current = (para.text || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n  para.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# This document contains a date/weekday heading paragraph followed by a\n# 20x5 table of simple arithmetic drills (one equation per cell). The edit\n# updates the heading date and rewrites every equation in the table.\n# Every \"old\" string below is unique within the document, so a whole-word,\n# case-sensitive Find/Replace targeting each one individually is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2023-02-23 Thursday', '2023-02-24 Friday'),\n    @('8+54=', '48-16='),\n    @('42+28=', '32+38='),\n    @('67-63=', '17+2='),\n    @('8+80=', '84-73='),\n    @('6+66=', '63+34='),\n    @('24+5=', '88-81='),\n    @('41-8=', '66-19='),\n    @('39+13=', '46-12='),\n    @('97-94=', '93-29='),\n    @('39+15=', '4+67='),\n    @('95-71=', '96-57='),\n    @('53-47=', '1+92='),\n    @('94-22=', '91-28='),\n    @('0+69=', '11+30='),\n    @('90-22=', '67+11='),\n    @('32+6=', '63-2='),\n    @('10+53=', '91-90='),\n    @('29+17=', '84-40='),\n    @('73-18=', '75-29='),\n    @('52+34=', '33+8='),\n    @('89-16=', '15+82='),\n    @('61-9=', '43+3='),\n    @('9+0=', '20+40='),\n    @('64-62=', '60+9='),\n    @('80-48=', '45-40='),\n    @('96-74=', '60-36='),\n    @('17+65=', '91-8='),\n    @('13+27=', '69-67='),\n    @('33+29=', '75+9='),\n    @('28+3=', '79-30='),\n    @('13+46=', '13+32='),\n    @('34-23=', '74+10='),\n    @('91-30=', '64-36='),\n    @('77-44=', '66-3='),\n    @('39+11=', '36-23='),\n    @('41-3=', '43-28='),\n    @('84-25=', '16+73='),\n    @('89-39=', '48-35='),\n    @('70-4=', '24+21='),\n    @('76-15=', '65+7='),\n    @('1+85=', '36-31='),\n    @('80-5=', '29+27='),\n    @('37+49=', '37+19='),\n    @('23+4=', '75+13='),\n    @('70-18=', '2+27='),\n    @('34+25=', '85-1='),\n    @('34+46=', '65+13='),\n    @('49+9=', '1+59='),\n    @('55-12=', '5+36='),\n    @('35-27=', '62-32='),\n    @('65+29=', '91-19='),\n    @('19+48=', '35+6='),\n    @('61-0=', '0+3='),\n    @('58+25=', '20+76='),\n    @('44+18=', '97-25='),\n    @('35-25=', '70-50='),\n    @('36+18=', '46-46='),\n    @('74-25=', '45+44='),\n    @('37+37=', '16+16='),\n    @('67-52=', '57+22='),\n    @('76-37=', '70+22='),\n    @('32+34=', '71+9='),\n    @('4+77=', '92-18='),\n    @('25+52=', '66-29='),\n    @('65-52=', '32+11='),\n    @('66+4=', '66+27='),\n    @('94-8=', '36+56='),\n    @('5+58=', '44+1='),\n    @('44+38=', '39+5='),\n    @('70-39=', '96-57='),\n    @('27-3=', '55+28='),\n    @('41+1=', '89-86='),\n    @('86-48=', '34+31='),\n    @('42+12=', '56+12='),\n    @('6+49=', '48-27='),\n    @('95-83=', '14+60='),\n    @('33+2=', '22+58='),\n    @('85-43=', '32+57='),\n    @('93-90=', '11+38='),\n    @('96-62=', '49-38='),\n    @('65-53=', '93-1='),\n    @('28+42=', '98-63='),\n    @('87-16=', '31+19='),\n    @('40+36=', '50+35='),\n    @('10+55=', '69-20='),\n    @('32+0=', '38+7='),\n    @('45-22=', '36-3='),\n    @('32+56=', '28+32='),\n    @('65+8=', '52-48='),\n    @('7+79=', '88-40='),\n    @('66+26=', '16+35='),\n    @('36-20=', '71+22='),\n    @('52-13=', '62-36='),\n    @('53+32=', '50+7='),\n    @('0+19=', '63-16='),\n    @('24-6=', '63+13='),\n    @('34+29=', '35+18='),\n    @('43-42=', '33+61='),\n    @('10+60=', '86-29='),\n    @('15+38=', '86-42=')\n)\n\n$successCount = 0\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceOne\n    if ($found) {\n        $successCount++\n    } else {\n        Write-Output \"WARNING: could not find text to replace: $oldText\"\n    }\n}\n\nWrite-Output \"Replaced $successCount of $($replacements.Count) items\"\n"}
